$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 1919.7059
$ws.Range("I86").Value = 2046.0834
$ws.Range("J86").Value = 1616.4
$ws.Range("K86").Value = 2046.0834
$ws.Range("L86").Value = 1616.4
$ws.Range("M86").Value = -923.0834
$ws.Range("N86").Value = -3862.4
# Row 89
$ws.Range("H89").Value = 1919.7059
$ws.Range("I89").Value = 2046.0834
$ws.Range("J89").Value = 1616.4
$ws.Range("K89").Value = 10230.417
$ws.Range("L89").Value = 8082
$ws.Range("M89").Value = -4614.416999999999
$ws.Range("N89").Value = -19314
# Row 125
$ws.Range("H125").Value = 5929.364
$ws.Range("I125").Value = 6253
$ws.Range("J125").Value = 5659.6665
$ws.Range("K125").Value = 56277
$ws.Range("L125").Value = 50936.9985
$ws.Range("M125").Value = -53817
$ws.Range("N125").Value = -55856.9985
# Row 134
$ws.Range("H134").Value = 26904.762
$ws.Range("J134").Value = 26904.762
$ws.Range("L134").Value = 26904.762
$ws.Range("N134").Value = -37044.762
# Row 138
$ws.Range("H138").Value = 4940.7046
$ws.Range("I138").Value = 1165.4117
$ws.Range("J138").Value = 5844.648
$ws.Range("K138").Value = 3496.2351
$ws.Range("L138").Value = 17533.944
$ws.Range("M138").Value = 1643.7649
$ws.Range("N138").Value = -27813.944

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6397.6343
$ws.Range("I32").Value = 5527.6
$ws.Range("K32").Value = 5527.6
$ws.Range("M32").Value = -5240.6
# Row 52
$ws.Range("H52").Value = 28999
$ws.Range("J52").Value = 28999
$ws.Range("L52").Value = 28999
$ws.Range("N52").Value = -29635
# Row 74
$ws.Range("H74").Value = 4886.2646
$ws.Range("I74").Value = 5389.7144
$ws.Range("J74").Value = 2536.8333
$ws.Range("K74").Value = 5389.7144
$ws.Range("L74").Value = 2536.8333
$ws.Range("M74").Value = -4515.7144
$ws.Range("N74").Value = -4284.8333
# Row 77
$ws.Range("H77").Value = 4886.2646
$ws.Range("I77").Value = 5389.7144
$ws.Range("J77").Value = 2536.8333
$ws.Range("K77").Value = 26948.572
$ws.Range("L77").Value = 12684.1665
$ws.Range("M77").Value = -22580.572
$ws.Range("N77").Value = -21420.1665
# Row 109
$ws.Range("H109").Value = 48800
$ws.Range("J109").Value = 48800
$ws.Range("L109").Value = 48800
$ws.Range("N109").Value = -51574
# Row 132
$ws.Range("H132").Value = 7626.9614
$ws.Range("I132").Value = 1751.3684
$ws.Range("J132").Value = 23575
$ws.Range("K132").Value = 5254.1052
$ws.Range("L132").Value = 70725
$ws.Range("M132").Value = -2724.1052
$ws.Range("N132").Value = -75785

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 52710
$ws.Range("J2").Value = 52710
$ws.Range("L2").Value = 52710
$ws.Range("N2").Value = -52936
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 500003100
$ws.Range("J25").Value = 500003100
$ws.Range("L25").Value = 500003100
$ws.Range("N25").Value = -500003448
# Row 31
$ws.Range("H31").Value = 17644.857
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 17644.857
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 17644.857
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -18234.857
# Row 34
$ws.Range("H34").Value = 17644.857
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 17644.857
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 17644.857
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -18048.857
# Row 62
$ws.Range("H62").Value = 8836.111000000001
$ws.Range("I62").Value = 7586.5
$ws.Range("J62").Value = 11335.333
$ws.Range("K62").Value = 7586.5
$ws.Range("L62").Value = 11335.333
$ws.Range("M62").Value = -6962.5
$ws.Range("N62").Value = -12583.333
# Row 65
$ws.Range("H65").Value = 8836.111000000001
$ws.Range("I65").Value = 7586.5
$ws.Range("J65").Value = 11335.333
$ws.Range("K65").Value = 37932.5
$ws.Range("L65").Value = 56676.665
$ws.Range("M65").Value = -34812.5
$ws.Range("N65").Value = -62916.665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Range("H112").Value = 2287.1365
$ws.Range("I112").Value = 472.33334
$ws.Range("J112").Value = 2573.6843
$ws.Range("K112").Value = 1417.00002
$ws.Range("L112").Value = 7721.0529
$ws.Range("M112").Value = -309.0000199999999
$ws.Range("N112").Value = -9937.052899999999
# Row 113
$ws.Range("H113").Value = 2857795
$ws.Range("I113").Value = 4167292.8
$ws.Range("J113").Value = 1111798.1
$ws.Range("K113").Value = 12501878.4
$ws.Range("L113").Value = 3335394.3
$ws.Range("M113").Value = -12499708.4
$ws.Range("N113").Value = -3339734.3
# Row 132
$ws.Range("H132").Value = 1822.6666
$ws.Range("I132").Value = 1655.1666
$ws.Range("J132").Value = 1857.9298
$ws.Range("K132").Value = 14896.4994
$ws.Range("L132").Value = 16721.3682
$ws.Range("M132").Value = -12366.4994
$ws.Range("N132").Value = -21781.3682

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 74
$ws.Range("H74").Value = 39999.59
$ws.Range("I74").Value = 39999
$ws.Range("K74").Value = 39999
$ws.Range("M74").Value = -39063
# Row 77
$ws.Range("H77").Value = 39999.59
$ws.Range("I77").Value = 39999
$ws.Range("K77").Value = 119997
$ws.Range("M77").Value = -115317
# Row 113
$ws.Range("H113").Value = 90910296
$ws.Range("I113").Value = 125001090
$ws.Range("J113").Value = 1530
$ws.Range("K113").Value = 125001090
$ws.Range("L113").Value = 1530
$ws.Range("M113").Value = -124998920
$ws.Range("N113").Value = -5870

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("K13").Value = 1000
$ws.Range("M13").Value = -860
# Row 46
$ws.Range("H46").Value = 10753654
$ws.Range("I46").Value = 22223048
$ws.Range("J46").Value = 1095.9375
$ws.Range("K46").Value = 22223048
$ws.Range("L46").Value = 1095.9375
$ws.Range("M46").Value = -22222860
$ws.Range("N46").Value = -1471.9375
# Row 61
$ws.Range("H61").Value = 2468.0715
$ws.Range("I61").Value = 2145.6667
$ws.Range("K61").Value = 2145.6667
$ws.Range("M61").Value = -1943.6667
# Row 113
$ws.Range("H113").Value = 2468.0715
$ws.Range("I113").Value = 2145.6667
$ws.Range("K113").Value = 2145.6667
$ws.Range("M113").Value = 24.33329999999978

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 47619524
$ws.Range("I107").Value = 66667084
$ws.Range("J107").Value = 624.8333
$ws.Range("K107").Value = 200001252
$ws.Range("L107").Value = 1874.4999
$ws.Range("M107").Value = -199999332
$ws.Range("N107").Value = -5714.4999
# Row 126
$ws.Range("H126").Value = 1550.3334
$ws.Range("I126").Value = 1139.4445
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 3418.3335
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -948.3335000000002
$ws.Range("N126").Value = -11440.0001
# Row 132
$ws.Range("H132").Value = 2386.25
$ws.Range("I132").Value = 2100.5
$ws.Range("J132").Value = 2481.5
$ws.Range("K132").Value = 6301.5
$ws.Range("L132").Value = 7444.5
$ws.Range("M132").Value = -3771.5
$ws.Range("N132").Value = -12504.5
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
